$wb = $excel.ActiveWorkbook

# --- Rename sheets -----------------------------------------------------
$wsMeasured  = $wb.Worksheets.Item("actual frequencies")
$wsMeasured.Name = "measured"

$wsPredicted = $wb.Worksheets.Item("predicted frequencies")
$wsPredicted.Name = "predicted"

$wsData = $wb.Worksheets.Item("data")

# --- Update header labels on "measured" (was "actual frequencies") -----
# Old headers were "Actual - short/medium/long" (shared strings 43-45);
# new headers are just "short"/"medium"/"long" (new shared strings).
$wsMeasured.Range("A1").Value = "short"
$wsMeasured.Range("B1").Value = "medium"
$wsMeasured.Range("C1").Value = "long"

# Remove the now-redundant "Frequency /Hz" label row and the predicted
# -frequency helper column (that data lives solely on "predicted" now).
$wsMeasured.Rows.Item(2).Delete()
$wsMeasured.Columns.Item(4).Delete()

# --- View / selection bookkeeping --------------------------------------
$wsMeasured.Range("C14").Select()

$wsData.Activate()
$wsData.Range("M9").Select()
$excel.ActiveWindow.ScrollColumn = 2
